$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (date / time number formats) of the last logged row (46)
# down onto the 6 new rows, so the new cells reuse the existing style
# records instead of creating duplicates.
$ws.Range("A46:E46").Copy()
$ws.Range("A47:E52").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 47: 11/01/2021, 17:00-18:00
$ws.Range("A47").Value = 44501
$ws.Range("B47").Value = 0.70833333333333337
$ws.Range("C47").Value = 0.75
$ws.Range("D47").Formula = "=C47-B47"

# Row 48: 11/03/2021, 09:12-09:57
$ws.Range("A48").Value = 44503
$ws.Range("B48").Value = 0.3833333333333333
$ws.Range("C48").Value = 0.4145833333333333
$ws.Range("D48").Formula = "=C48-B48"

# Row 49: 11/03/2021, 18:15-19:00
$ws.Range("A49").Value = 44503
$ws.Range("B49").Value = 0.76041666666666663
$ws.Range("C49").Value = 0.79166666666666663
$ws.Range("D49").Formula = "=C49-B49"

# Row 50: 11/04/2021, 09:12-09:57
$ws.Range("A50").Value = 44504
$ws.Range("B50").Value = 0.3833333333333333
$ws.Range("C50").Value = 0.4145833333333333
$ws.Range("D50").Formula = "=C50-B50"

# Row 51: 11/05/2021, 09:12-09:57
$ws.Range("A51").Value = 44505
$ws.Range("B51").Value = 0.3833333333333333
$ws.Range("C51").Value = 0.4145833333333333
$ws.Range("D51").Formula = "=C51-B51"

# Row 52: 11/07/2021, 12:00-14:00
$ws.Range("A52").Value = 44507
$ws.Range("B52").Value = 0.5
$ws.Range("C52").Value = 0.58333333333333337
$ws.Range("D52").Formula = "=C52-B52"

# Rows 48-51 (the presentation-related entries) get their start/end times
# recolored to plain black instead of the theme color.
$ws.Range("B48:C48").Font.Color = 0
$ws.Range("B49:C49").Font.Color = 0
$ws.Range("B50:C50").Font.Color = 0
$ws.Range("B51:C51").Font.Color = 0

# Activity notes, in the order they were filled in.
$ws.Range("E48").Value = "marcus's presentation"
$ws.Range("E50").Value = "Jimena's presentation"
$ws.Range("E51").Value = "Sarah's presentation"
$ws.Range("E47").Value = "Emailed mentor"
$ws.Range("E52").Value = "worked on STS application, helped sean with STS application"
$ws.Range("E49").Value = "worked on presentation"

# The biweekly-total formula at the top now sums over the new block of rows.
$ws.Range("J1").Formula = "=SUM(D41:D55)"

# Scroll the view down to where the new entries were added and leave the
# selection on the last-edited cell.
$excel.ActiveWindow.ScrollRow = 40
$ws.Range("D50").Select()

Write-Output "edit applied"
